$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- B2 / C2 become numeric 0 instead of text "-" ---
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# --- New Area column (G) segment-area formulas ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4").Formula = "=(D4-D3)*B4/100"
$ws.Range("G5").Formula = "=(D5-D4)*B5/100"
$ws.Range("G6").Formula = "=(D6-D5)*B6/100"
$ws.Range("G7").Formula = "=(D7-D6)*B7/100"
$ws.Range("G8").Formula = "=(D8-D7)*B8/100"
$ws.Range("G9").Formula = "=(D9-D8)*B9/100"
$ws.Range("G10").Formula = "=(D10-D9)*B10/100"
$ws.Range("G11").Formula = "=(D11-D10)*B11/100"
$ws.Range("G12").Formula = "=(D12-D11)*B12/100"
$ws.Range("G13").Formula = "=(D13-D12)*B13/100"
$ws.Range("G14").Formula = "=(D14-D13)*B14/100"
$ws.Range("G15").Formula = "=(D15-D14)*B15/100"

# --- Total area (H2) and summary mirrors (J2/K2) ---
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- E4 no longer shares E3's formula, gets its own explicit formula ---
$ws.Range("E4").Formula = "=(D4-D3)*(B4/100)*C4"

# --- Update selection to match the saved view state ---
$ws.Range("E4").Select()
